$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before A ("Match ID"), shifting every existing column
# (A..W) one slot to the right (B..X). Excel's native Insert keeps per-cell
# styles/values attached to the shifted cells, which is what we want here.
$ws.Columns("A").Insert()

# Header cell for the new column.
$ws.Range("A1").Value = "Match ID"

# New column A uses the bold font (no border/alignment) - same font as the
# existing header style, just without the border/alignment tweaks.
$ws.Range("A1:A19").Font.Bold = $true

# Row 2 is a hidden "field name" helper row; the new column stays blank there.
# Toggle Hidden off/on around the write so the engine doesn't recompute an
# explicit row height for a hidden row (which would otherwise stamp a
# ht="9.5" customHeight="1" we don't want).
$ws.Rows.Item(2).Hidden = $false
$ws.Cells.Item(2, 1).Value = ""
$ws.Rows.Item(2).Hidden = $true

# Row 3 is a hidden spacer row that previously had no cells at all; give it
# a blank, styled A3 cell to match.
$ws.Rows.Item(3).Hidden = $false
$ws.Cells.Item(3, 1).Value = ""
$ws.Rows.Item(3).Hidden = $true

# Data rows 4-19: the match id value.
for ($r = 4; $r -le 19; $r++) {
    $ws.Cells.Item($r, 1).Value = 5
}

# Row 20 is the hidden totals row; it gets the match id too, but keeps the
# default (non-bold) style, unlike rows 1-19.
$ws.Rows.Item(20).Hidden = $false
$ws.Cells.Item(20, 1).Value = 5
$ws.Cells.Item(20, 1).Font.Bold = $false
$ws.Rows.Item(20).Hidden = $true

# Match the saved selection/active cell.
$ws.Range("F24").Select() | Out-Null
